$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.287.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.03%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.369.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.02%  "

# Row 4
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.46%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("E8").Value = "  +5.50%  "

# Row 9
$ws.Range("E9").Value = "  +4.39%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.31%  "

# Row 11
$ws.Range("E11").Value = "  -1.95%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.81%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.89%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.790.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.99%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.227.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.83%  "

# Row 16
$ws.Range("E16").Value = "  +2.43%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.351.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.65%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.57%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.93%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "331.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.68%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.58%  "

# Row 22
$ws.Range("E22").Value = "  +0.02%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "63.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.15%  "

# Row 24
$ws.Range("E24").Value = "  -0.84%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.16%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.49%  "

# Row 27
$ws.Range("E27").Value = "  -6.28%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.10%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.22%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0744"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.90%  "

# Row 31
$ws.Range("E31").Value = "  +0.39%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.20%  "

# Row 33
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.01%  "

# Row 34
$ws.Range("B34").Value = "SuiNetwork"
$ws.Range("C34").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.19%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.53%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.13%  "

# Row 37
$ws.Range("E37").Value = "  -1.83%  "

# Row 38
$ws.Range("E38").Value = "  -1.56%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.415"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.26%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "142.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.69%  "

# Row 41
$ws.Range("E41").Value = "  +2.16%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "287.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.99%  "

# Row 43
$ws.Range("E43").Value = "  +2.98%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0520"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.93%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.18%  "

# Row 46
$ws.Range("E46").Value = "  +0.92%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0223"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.77%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.390"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.09%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.28%  "

# Row 50
$ws.Range("E50").Value = "  +0.71%  "

# Row 51
$ws.Range("E51").Value = "  +0.46%  "
